$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.155.87"
$ws.Range("E2").Value = "  -2.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.868.22"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.13"
$ws.Range("E5").Value = "  -1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5170"
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3759"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07152"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8917"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.76"
$ws.Range("E11").Value = "  -2.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07537"
$ws.Range("E12").Value = "  -1.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.867.52"
$ws.Range("E13").Value = "  -1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.305"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.58"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9998"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008483"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.08"
$ws.Range("E18").Value = "  -3.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9998"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.186.68"
$ws.Range("E20").Value = "  -2.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.996"
$ws.Range("E21").Value = "  -2.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.101.40"
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.460"
$ws.Range("E24").Value = "  -2.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.839"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.89"
$ws.Range("E26").Value = "  -5.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.95"
$ws.Range("E27").Value = "  -1.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.084"
$ws.Range("E28").Value = "  -3.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.92"
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.662"
$ws.Range("E30").Value = "  -3.93%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.682"
$ws.Range("E31").Value = "  -3.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09248"
$ws.Range("E32").Value = "  +2.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05127"
$ws.Range("E33").Value = "  -2.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.082"
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.162"
$ws.Range("E35").Value = "  -5.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7254"
$ws.Range("E36").Value = "  -7.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02032"
$ws.Range("E37").Value = "  -2.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.104"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.504"
$ws.Range("E39").Value = "  -3.53%  "
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5286"
$ws.Range("E41").Value = "  -5.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.504"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.68"
$ws.Range("E43").Value = "  +1.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.310"
$ws.Range("E44").Value = "  -2.46%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1470"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9995"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4616"
$ws.Range("E47").Value = "  -4.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.972"
$ws.Range("E48").Value = "  -4.98%  "
$ws.Range("E49").Value = "  -3.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.69"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.62"
